$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300.33334
$ws.Range("I2").Value = 199.5
$ws.Range("J2").Value = 502
$ws.Range("K2").Value = 199.5
$ws.Range("L2").Value = 502
$ws.Range("M2").Value = -86.5
$ws.Range("N2").Value = -728
$ws.Range("H18").Value = 987.1667
$ws.Range("I18").Value = 710.6667
$ws.Range("K18").Value = 710.6667
$ws.Range("M18").Value = -426.6667
$ws.Range("H38").Value = 1754504
$ws.Range("I38").Value = 2222305
$ws.Range("J38").Value = 250
$ws.Range("K38").Value = 6666915
$ws.Range("L38").Value = 750
$ws.Range("M38").Value = -6666543
$ws.Range("N38").Value = -1494
$ws.Range("H40").Value = 2207.1428
$ws.Range("I40").Value = 1499.8572
$ws.Range("K40").Value = 1499.8572
$ws.Range("M40").Value = -1324.8572
$ws.Range("H58").Value = 1227095
$ws.Range("I58").Value = 2451240
$ws.Range("J58").Value = 2950
$ws.Range("K58").Value = 7353720
$ws.Range("L58").Value = 8850
$ws.Range("M58").Value = -7353570
$ws.Range("N58").Value = -9150
$ws.Range("H62").Value = 2947.3684
$ws.Range("I62").Value = 2866.2666
$ws.Range("K62").Value = 2866.2666
$ws.Range("M62").Value = -2242.2666
$ws.Range("H64").Value = 69626.664
$ws.Range("I64").Value = 113777.78
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 113777.78
$ws.Range("L64").Value = 3400
$ws.Range("M64").Value = -113529.78
$ws.Range("N64").Value = -3896
$ws.Range("H65").Value = 2947.3684
$ws.Range("I65").Value = 2866.2666
$ws.Range("K65").Value = 14331.333
$ws.Range("M65").Value = -11211.333
$ws.Range("H67").Value = 69626.664
$ws.Range("I67").Value = 113777.78
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 113777.78
$ws.Range("L67").Value = 3400
$ws.Range("M67").Value = -112919.78
$ws.Range("N67").Value = -5116
$ws.Range("H107").Value = 708.05
$ws.Range("I107").Value = 650.5789
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 650.5789
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 1269.4211
$ws.Range("N107").Value = -5640
$ws.Range("H137").Value = 1622.0555
$ws.Range("I137").Value = 1542
$ws.Range("J137").Value = 1902.25
$ws.Range("K137").Value = 4626
$ws.Range("L137").Value = 5706.75
$ws.Range("M137").Value = -2076
$ws.Range("N137").Value = -10806.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11772.225
$ws.Range("I32").Value = 10158.779
$ws.Range("J32").Value = 23335.25
$ws.Range("K32").Value = 10158.779
$ws.Range("L32").Value = 23335.25
$ws.Range("M32").Value = -9871.779
$ws.Range("N32").Value = -23909.25
$ws.Range("H45").Value = 72875.36
$ws.Range("I45").Value = 167568.5
$ws.Range("J45").Value = 1855.5
$ws.Range("K45").Value = 167568.5
$ws.Range("L45").Value = 1855.5
$ws.Range("M45").Value = -167191.5
$ws.Range("N45").Value = -2609.5
$ws.Range("H61").Value = 1723.6765
$ws.Range("I61").Value = 1481.15
$ws.Range("J61").Value = 2070.1428
$ws.Range("K61").Value = 1481.15
$ws.Range("L61").Value = 2070.1428
$ws.Range("M61").Value = -1269.15
$ws.Range("N61").Value = -2494.1428
$ws.Range("H64").Value = 43980.332
$ws.Range("J64").Value = 43980.332
$ws.Range("L64").Value = 43980.332
$ws.Range("N64").Value = -44476.332
$ws.Range("H67").Value = 43980.332
$ws.Range("J67").Value = 43980.332
$ws.Range("L67").Value = 43980.332
$ws.Range("N67").Value = -45696.332
$ws.Range("H74").Value = 1558.2858
$ws.Range("I74").Value = 1513.5294
$ws.Range("K74").Value = 1513.5294
$ws.Range("M74").Value = -639.5293999999999
$ws.Range("H77").Value = 1558.2858
$ws.Range("I77").Value = 1513.5294
$ws.Range("K77").Value = 7567.646999999999
$ws.Range("M77").Value = -3199.646999999999
$ws.Range("H88").Value = 1699.6
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 1749.5
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1749.5
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2561.5
$ws.Range("H91").Value = 1699.6
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 1749.5
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1749.5
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4557.5
$ws.Range("H122").Value = 2811.25
$ws.Range("I122").Value = 2517.2222
$ws.Range("K122").Value = 7551.6666
$ws.Range("M122").Value = -5101.6666
$ws.Range("H136").Value = 1723.6765
$ws.Range("I136").Value = 1481.15
$ws.Range("J136").Value = 2070.1428
$ws.Range("K136").Value = 4443.450000000001
$ws.Range("L136").Value = 6210.428400000001
$ws.Range("M136").Value = -1893.450000000001
$ws.Range("N136").Value = -11310.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 68000
$ws.Range("J57").Value = 68000
$ws.Range("L57").Value = 68000
$ws.Range("N57").Value = -69440
$ws.Range("H86").Value = 187583.33
$ws.Range("I86").Value = 224740
$ws.Range("K86").Value = 224740
$ws.Range("M86").Value = -223617
$ws.Range("H89").Value = 187583.33
$ws.Range("I89").Value = 224740
$ws.Range("K89").Value = 1123700
$ws.Range("M89").Value = -1118084
$ws.Range("H107").Value = 166746850
$ws.Range("I107").Value = 250114910
$ws.Range("J107").Value = 10749.5
$ws.Range("K107").Value = 250114910
$ws.Range("L107").Value = 10749.5
$ws.Range("M107").Value = -250112990
$ws.Range("N107").Value = -14589.5
$ws.Range("H135").Value = 68000
$ws.Range("J135").Value = 68000
$ws.Range("L135").Value = 68000
$ws.Range("N135").Value = -78140
$ws.Range("H136").Value = 68000
$ws.Range("J136").Value = 68000
$ws.Range("L136").Value = 68000
$ws.Range("N136").Value = -78200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 716.1667
$ws.Range("J16").Value = 756.5
$ws.Range("L16").Value = 756.5
$ws.Range("N16").Value = -1330.5
$ws.Range("H31").Value = 30587.549
$ws.Range("I31").Value = 1301.1714
$ws.Range("K31").Value = 1301.1714
$ws.Range("M31").Value = -1006.1714
$ws.Range("H34").Value = 30587.549
$ws.Range("I34").Value = 1301.1714
$ws.Range("K34").Value = 1301.1714
$ws.Range("M34").Value = -1099.1714
$ws.Range("H41").Value = 10300
$ws.Range("J41").Value = 14850
$ws.Range("L41").Value = 14850
$ws.Range("N41").Value = -15706
$ws.Range("H50").Value = 9457.6
$ws.Range("J50").Value = 9457.6
$ws.Range("L50").Value = 9457.6
$ws.Range("N50").Value = -10707.6
$ws.Range("H51").Value = 8099
$ws.Range("J51").Value = 8099
$ws.Range("L51").Value = 8099
$ws.Range("N51").Value = -9571
$ws.Range("H60").Value = 15020.75
$ws.Range("J60").Value = 15020.75
$ws.Range("L60").Value = 15020.75
$ws.Range("N60").Value = -16042.75
$ws.Range("H61").Value = 8099
$ws.Range("J61").Value = 8099
$ws.Range("L61").Value = 8099
$ws.Range("N61").Value = -8795
$ws.Range("H68").Value = 17187.422
$ws.Range("J68").Value = 17187.422
$ws.Range("L68").Value = 17187.422
$ws.Range("N68").Value = -18685.422
$ws.Range("H71").Value = 17187.422
$ws.Range("J71").Value = 17187.422
$ws.Range("L71").Value = 51562.266
$ws.Range("N71").Value = -59050.266
$ws.Range("H74").Value = 40422
$ws.Range("J74").Value = 40422
$ws.Range("L74").Value = 40422
$ws.Range("N74").Value = -42170
$ws.Range("H77").Value = 40422
$ws.Range("J77").Value = 40422
$ws.Range("L77").Value = 121266
$ws.Range("N77").Value = -130002
$ws.Range("H99").Value = 2313.0688
$ws.Range("I99").Value = 1784.2858
$ws.Range("K99").Value = 1784.2858
$ws.Range("M99").Value = -286.2858000000001
$ws.Range("H113").Value = 716.1667
$ws.Range("J113").Value = 756.5
$ws.Range("L113").Value = 756.5
$ws.Range("N113").Value = -5096.5
$ws.Range("H122").Value = 1147
$ws.Range("H126").Value = 2313.0688
$ws.Range("I126").Value = 1784.2858
$ws.Range("K126").Value = 5352.857400000001
$ws.Range("M126").Value = -2882.857400000001
$ws.Range("H140").Value = 66500
$ws.Range("J140").Value = 66500
$ws.Range("L140").Value = 66500
$ws.Range("N140").Value = -76860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1287.2106
$ws.Range("I107").Value = 605.1818
$ws.Range("J107").Value = 2225
$ws.Range("K107").Value = 1815.5454
$ws.Range("L107").Value = 6675
$ws.Range("M107").Value = 104.4546
$ws.Range("N107").Value = -10515
$ws.Range("H113").Value = 1251.4
$ws.Range("I113").Value = 2524
$ws.Range("J113").Value = 615.1
$ws.Range("K113").Value = 7572
$ws.Range("L113").Value = 1845.3
$ws.Range("M113").Value = -5402
$ws.Range("N113").Value = -6185.3
$ws.Range("H124").Value = 3447.5
$ws.Range("I124").Value = 1915
$ws.Range("J124").Value = 4980
$ws.Range("K124").Value = 5745
$ws.Range("L124").Value = 14940
$ws.Range("M124").Value = -835
$ws.Range("N124").Value = -24760
$ws.Range("H125").Value = 1400
$ws.Range("I125").Value = 1100
$ws.Range("K125").Value = 3300
$ws.Range("M125").Value = 1620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2836.3635
$ws.Range("I122").Value = 2537.5
$ws.Range("K122").Value = 7612.5
$ws.Range("M122").Value = -5162.5
$ws.Range("H135").Value = 44850.117
$ws.Range("J135").Value = 44850.117
$ws.Range("L135").Value = 44850.117
$ws.Range("N135").Value = -54990.117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2634.1428
$ws.Range("I7").Value = 1796.8
$ws.Range("J7").Value = 3099.3333
$ws.Range("K7").Value = 1796.8
$ws.Range("L7").Value = 3099.3333
$ws.Range("M7").Value = -1684.8
$ws.Range("N7").Value = -3323.3333
$ws.Range("H46").Value = 4674.3335
$ws.Range("I46").Value = 5546.6665
$ws.Range("J46").Value = 4092.7778
$ws.Range("K46").Value = 5546.6665
$ws.Range("L46").Value = 4092.7778
$ws.Range("M46").Value = -5358.6665
$ws.Range("N46").Value = -4468.7778
$ws.Range("H122").Value = 5848.4287
$ws.Range("J122").Value = 6187
$ws.Range("L122").Value = 18561
$ws.Range("N122").Value = -23461
$ws.Range("H126").Value = 2634.1428
$ws.Range("I126").Value = 1796.8
$ws.Range("J126").Value = 3099.3333
$ws.Range("K126").Value = 5390.4
$ws.Range("L126").Value = 9297.999899999999
$ws.Range("M126").Value = -2920.4
$ws.Range("N126").Value = -14237.9999
$ws.Range("H132").Value = 4971.1904
$ws.Range("I132").Value = 7990.6665
$ws.Range("J132").Value = 2706.5833
$ws.Range("K132").Value = 23971.9995
$ws.Range("L132").Value = 8119.749899999999
$ws.Range("M132").Value = -21441.9995
$ws.Range("N132").Value = -13179.7499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7695417
$ws.Range("I62").Value = 12823512
$ws.Range("K62").Value = 12823512
$ws.Range("M62").Value = -12822888
$ws.Range("H65").Value = 7695417
$ws.Range("I65").Value = 12823512
$ws.Range("K65").Value = 64117560
$ws.Range("M65").Value = -64114440
$ws.Range("H126").Value = 1605.4286
$ws.Range("I126").Value = 1706.3334
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 5119.0002
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -2649.0002
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 2585
$ws.Range("I132").Value = 2771.8215
$ws.Range("J132").Value = 2182.6155
$ws.Range("K132").Value = 8315.4645
$ws.Range("L132").Value = 6547.8465
$ws.Range("M132").Value = -5785.4645
$ws.Range("N132").Value = -11607.8465
$ws.Range("H136").Value = 1538.2373
$ws.Range("I136").Value = 715.56665
$ws.Range("J136").Value = 2389.276
$ws.Range("K136").Value = 2146.69995
$ws.Range("L136").Value = 7167.828
$ws.Range("M136").Value = 403.3000499999998
$ws.Range("N136").Value = -12267.828
